# AMOVA - isolation by distance - allele frequencies
# Add the missing Hannover coordinate row (lat/long/elevation) on Tabelle1
# and move the selection to A26, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C24").Value = 52.3759
$ws.Range("D24").Value = 9.732
$ws.Range("E24").Value = 55

$ws.Range("A26").Select()
